$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 7")

# --- Resize / reposition the "TextBox 7" members list box ---
$shp.Left   = 5638800 / 12700
$shp.Top    = 3540711 / 12700
$shp.Width  = 3370603 / 12700
$shp.Height = 2031325 / 12700

# --- Insert " Trưởng" (red) right after "Đội" so "(Tú Đội)" becomes
#     "(Tú Đội Trưởng)" ---
$tr = $shp.TextFrame.TextRange
$doi = $tr.Find("Đội")
$afterDoi = $tr.Characters($doi.Start + $doi.Length, 1)
$null = $afterDoi.InsertBefore(" Trưởng")

# The inserted text inherits the red formatting of its neighbours already,
# but make it explicit so the new runs are unambiguously red (FF0000),
# matching "(Tú Đội)" styling.
$spaceRun = $tr.Characters($doi.Start + $doi.Length, 1)
$spaceRun.Font.Color.RGB = 255

$wordRun = $tr.Characters($doi.Start + $doi.Length + 1, 6)
$wordRun.Font.Color.RGB = 255
